$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Column width fix for columns I (9) and J (10) to match column H (8) ---
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# --- Crime-statistics table data updates (rows 15-28) ---
# Row 15
$ws.Range("M15").Value = 100

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 39
$ws.Range("K16").Value = -31.578947368421
$ws.Range("L16").Value = -29.090909090909
$ws.Range("N16").Value = -87.213114754098

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 1
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("E17").Value = 300
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 118
$ws.Range("J17").Value = 85
$ws.Range("K17").Value = 38.823529411764
$ws.Range("L17").Value = 29.670329670329
$ws.Range("M17").Value = 145.833333333333
$ws.Range("N17").Value = -6.349206349206

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -30
$ws.Range("F18").Value = 34
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = 9.677419354838
$ws.Range("I18").Value = 299
$ws.Range("J18").Value = 288
$ws.Range("K18").Value = 3.819444444444
$ws.Range("L18").Value = 0.335570469798
$ws.Range("M18").Value = 29.437229437229
$ws.Range("N18").Value = -67.214912280701

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 92
$ws.Range("I19").Value = 472
$ws.Range("J19").Value = 419
$ws.Range("K19").Value = 12.649164677804
$ws.Range("L19").Value = -17.62652705061
$ws.Range("M19").Value = 38.823529411764
$ws.Range("N19").Value = -9.923664122137

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 125
$ws.Range("I20").Value = 289
$ws.Range("J20").Value = 233
$ws.Range("K20").Value = 24.034334763948
$ws.Range("L20").Value = 61.45251396648
$ws.Range("M20").Value = 127.55905511811
$ws.Range("N20").Value = -90.509031198686

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 44.444444444444
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 76
$ws.Range("H21").Value = 56.578947368421
$ws.Range("I21").Value = 1231
$ws.Range("J21").Value = 1091
$ws.Range("K21").Value = 12.832263978001
$ws.Range("L21").Value = 1.9884009942
$ws.Range("M21").Value = 48.313253012048
$ws.Range("N21").Value = -75.015222244773

# Row 24
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 505
$ws.Range("J24").Value = 518
$ws.Range("K24").Value = -2.509652509652
$ws.Range("L24").Value = -5.783582089552
$ws.Range("M24").Value = 22.27602905569

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -61.111111111111
$ws.Range("I25").Value = 93
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = -23.140495867768
$ws.Range("L25").Value = -16.964285714285

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 189
$ws.Range("J26").Value = 206
$ws.Range("K26").Value = -8.252427184466
$ws.Range("L26").Value = -10.849056603773
$ws.Range("M26").Value = 18.125

# Row 28
$ws.Range("G28").Value = 1

